$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 21742070
$ws.Range("I32").Value = 83336740
$ws.Range("K32").Value = 83336740
$ws.Range("M32").Value = -83336414

$ws.Range("H51").Value = 5769.72
$ws.Range("J51").Value = 8424.299999999999
$ws.Range("L51").Value = 8424.299999999999
$ws.Range("N51").Value = -9392.299999999999

$ws.Range("H74").Value = 4666.6665
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 4666.6665
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H98").Value = 2080.55
$ws.Range("I98").Value = 1876.5625
$ws.Range("K98").Value = 1876.5625
$ws.Range("M98").Value = -378.5625

$ws.Range("H103").Value = 1440.8
$ws.Range("I103").Value = 1200
$ws.Range("J103").Value = 1501
$ws.Range("K103").Value = 3600
$ws.Range("L103").Value = 4503
$ws.Range("M103").Value = -3014
$ws.Range("N103").Value = -5675

$ws.Range("H107").Value = 582.1667
$ws.Range("I107").Value = 598.6
$ws.Range("K107").Value = 598.6
$ws.Range("M107").Value = 1321.4

$ws.Range("H122").Value = 2080.55
$ws.Range("I122").Value = 1876.5625
$ws.Range("K122").Value = 5629.6875
$ws.Range("M122").Value = -3179.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 643194.8
$ws.Range("I2").Value = 1028801.1
$ws.Range("J2").Value = 47257.727
$ws.Range("K2").Value = 1028801.1
$ws.Range("L2").Value = 47257.727
$ws.Range("M2").Value = -1028688.1
$ws.Range("N2").Value = -47483.727

$ws.Range("H63").Value = 3335
$ws.Range("I63").Value = 2752.5
$ws.Range("K63").Value = 2752.5
$ws.Range("M63").Value = -2066.5

$ws.Range("H66").Value = 3335
$ws.Range("I66").Value = 2752.5
$ws.Range("K66").Value = 13762.5
$ws.Range("M66").Value = -10330.5

$ws.Range("H116").Value = 643194.8
$ws.Range("I116").Value = 1028801.1
$ws.Range("J116").Value = 47257.727
$ws.Range("K116").Value = 1028801.1
$ws.Range("L116").Value = 47257.727
$ws.Range("M116").Value = -1026507.1
$ws.Range("N116").Value = -51845.727

$ws.Range("H122").Value = 3834.8064
$ws.Range("I122").Value = 2051.875
$ws.Range("K122").Value = 6155.625
$ws.Range("M122").Value = -3705.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 643194.8
$ws.Range("I3").Value = 1028801.1
$ws.Range("J3").Value = 47257.727
$ws.Range("K3").Value = 1028801.1
$ws.Range("L3").Value = 47257.727
$ws.Range("M3").Value = -1028687.1
$ws.Range("N3").Value = -47485.727

$ws.Range("H20").Value = 2990.111
$ws.Range("I20").Value = 2613.875
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 2613.875
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -2366.875
$ws.Range("N20").Value = -6494

$ws.Range("H86").Value = 3361.5173
$ws.Range("J86").Value = 3629.1428
$ws.Range("L86").Value = 3629.1428
$ws.Range("N86").Value = -5875.1428

$ws.Range("H89").Value = 3361.5173
$ws.Range("J89").Value = 3629.1428
$ws.Range("L89").Value = 18145.714
$ws.Range("N89").Value = -29377.714

$ws.Range("H134").Value = 3428.5789
$ws.Range("I134").Value = 2909.6667
$ws.Range("K134").Value = 8729.000100000001
$ws.Range("M134").Value = -6194.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1919.6428
$ws.Range("I58").Value = 2056.8572
$ws.Range("J58").Value = 1782.4286
$ws.Range("K58").Value = 2056.8572
$ws.Range("L58").Value = 1782.4286
$ws.Range("M58").Value = -1853.8572
$ws.Range("N58").Value = -2188.4286

$ws.Range("H105").Value = 1749167
$ws.Range("I105").Value = 2273417
$ws.Range("J105").Value = 1666.6666
$ws.Range("K105").Value = 2273417
$ws.Range("L105").Value = 1666.6666
$ws.Range("M105").Value = -2271670
$ws.Range("N105").Value = -5160.6666

$ws.Range("H107").Value = 675735.4399999999
$ws.Range("I107").Value = 1299722.1
$ws.Range("K107").Value = 1299722.1
$ws.Range("M107").Value = -1297802.1

$ws.Range("H122").Value = 4547.3335
$ws.Range("I122").Value = 2689.4285
$ws.Range("J122").Value = 6173
$ws.Range("K122").Value = 8068.2855
$ws.Range("L122").Value = 18519
$ws.Range("M122").Value = -5618.2855
$ws.Range("N122").Value = -23419

$ws.Range("H136").Value = 1919.6428
$ws.Range("I136").Value = 2056.8572
$ws.Range("J136").Value = 1782.4286
$ws.Range("K136").Value = 6170.571599999999
$ws.Range("L136").Value = 5347.2858
$ws.Range("M136").Value = -3620.571599999999
$ws.Range("N136").Value = -10447.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 276.9091
$ws.Range("I33").Value = 257
$ws.Range("J33").Value = 311.75
$ws.Range("K33").Value = 1542
$ws.Range("L33").Value = 1870.5
$ws.Range("M33").Value = -1259
$ws.Range("N33").Value = -2436.5

$ws.Range("H57").Value = 8354.6
$ws.Range("I57").Value = 4193.5
$ws.Range("K57").Value = 12580.5
$ws.Range("M57").Value = -12021.5

$ws.Range("H121").Value = 687.7778
$ws.Range("J121").Value = 417.4
$ws.Range("L121").Value = 1252.2
$ws.Range("N121").Value = -3872.2

$ws.Range("H132").Value = 2938.8572
$ws.Range("J132").Value = 4390
$ws.Range("L132").Value = 39510
$ws.Range("N132").Value = -44570

$ws.Range("H136").Value = 774819.6
$ws.Range("I136").Value = 911150.4399999999
$ws.Range("K136").Value = 2733451.32
$ws.Range("M136").Value = -2728351.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 50000.668
$ws.Range("J117").Value = 50000.668
$ws.Range("L117").Value = 50000.668
$ws.Range("N117").Value = -56884.668

$ws.Range("H126").Value = 3301.8647
$ws.Range("I126").Value = 2133.8333
$ws.Range("J126").Value = 5458.231
$ws.Range("K126").Value = 6401.499899999999
$ws.Range("L126").Value = 16374.693
$ws.Range("M126").Value = -3931.499899999999
$ws.Range("N126").Value = -21314.693

$ws.Range("H132").Value = 8118.875
$ws.Range("I132").Value = 7594.6
$ws.Range("J132").Value = 8992.666999999999
$ws.Range("K132").Value = 22783.8
$ws.Range("L132").Value = 26978.001
$ws.Range("M132").Value = -20253.8
$ws.Range("N132").Value = -32038.001

$ws.Range("H138").Value = 59849
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5799.4287
$ws.Range("I40").Value = 4421.25
$ws.Range("K40").Value = 4421.25
$ws.Range("M40").Value = -4285.25

$ws.Range("H59").Value = 54400
$ws.Range("J59").Value = 54400
$ws.Range("L59").Value = 54400
$ws.Range("N59").Value = -55708

$ws.Range("H111").Value = 91887
$ws.Range("J111").Value = 91887
$ws.Range("L111").Value = 91887
$ws.Range("N111").Value = -100067

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28569.75
$ws.Range("J41").Value = 28569.75
$ws.Range("L41").Value = 28569.75
$ws.Range("N41").Value = -29349.75

$ws.Range("H81").Value = 4172492.5
$ws.Range("I81").Value = 3478071.2
$ws.Range("J81").Value = 5214124.5
$ws.Range("K81").Value = 6956142.4
$ws.Range("L81").Value = 10428249
$ws.Range("M81").Value = -6955081.4
$ws.Range("N81").Value = -10430371

$ws.Range("H84").Value = 4172492.5
$ws.Range("I84").Value = 3478071.2
$ws.Range("J84").Value = 5214124.5
$ws.Range("K84").Value = 34780712
$ws.Range("L84").Value = 52141245
$ws.Range("M84").Value = -34775408
$ws.Range("N84").Value = -52151853

$ws.Range("H107").Value = 2059.2
$ws.Range("I107").Value = 2168.4443
$ws.Range("J107").Value = 1895.3334
$ws.Range("K107").Value = 6505.3329
$ws.Range("L107").Value = 5686.0002
$ws.Range("M107").Value = -4585.3329
$ws.Range("N107").Value = -9526.0002

$ws.Range("H116").Value = 57154.8
$ws.Range("J116").Value = 57154.8
$ws.Range("L116").Value = 57154.8
$ws.Range("N116").Value = -66332.8

$ws.Range("H123").Value = 92429
$ws.Range("J123").Value = 92429
$ws.Range("L123").Value = 92429
$ws.Range("N123").Value = -102229

$ws.Range("H141").Value = 47381.332
$ws.Range("J141").Value = 47381.332
$ws.Range("L141").Value = 47381.332
$ws.Range("N141").Value = -57741.332
